# Trade #32 closed at 2026-02-17 12:38:31 - unknown UNKNOWN +0.000%
#
# This script updates the "live_trading_results" workbook to record the
# closing of trade #32 on the MarketMaking strategy:
#   - Summary sheet metrics are refreshed
#   - Strategy Status row for MarketMaking is refreshed
#   - A new trade row (#32 / row 33) is appended to both the
#     "All Trades" sheet and the "MarketMaking" sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.72   # Current Capital
$summary.Range("B4").Value = 0.71      # Total P&L $
$summary.Range("B5").Value = 0.44      # Total P&L %
$summary.Range("B6").Value = 32        # Total Trades
$summary.Range("B8").Value = 11        # Losing Trades
$summary.Range("B9").Value = 40.62     # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet (MarketMaking row, row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.72     # Capital
$status.Range("D4").Value = 32         # Trades
$status.Range("E4").Value = 0.71       # P&L $
$status.Range("F4").Value = 0.72       # P&L %
$status.Range("G4").Value = 40.62      # Win Rate %

# ---------------------------------------------------------------------
# 3. Append the new closed trade (#32) to "All Trades" and
#    "MarketMaking" sheets as row 33
# ---------------------------------------------------------------------
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A33").Value = 32

    # Force the date-looking text to remain plain text instead of being
    # auto-converted into a date serial number.
    $ws.Range("B33").NumberFormat = "@"
    $ws.Range("B33").Value = "2026-02-17"

    $ws.Range("C33").Value = "12:38:25"
    $ws.Range("D33").Value = "MarketMaking"
    $ws.Range("E33").Value = "DOWN"
    $ws.Range("F33").Value = 0.31
    $ws.Range("G33").Value = 0.253884
    $ws.Range("H33").Value = "CLOSED"
    $ws.Range("I33").Value = -18.1018
    $ws.Range("J33").Value = -0.06
    $ws.Range("K33").Value = 100.72
    $ws.Range("L33").Value = 0
    $ws.Range("M33").Value = 0
    $ws.Range("N33").Value = 0.6
    $ws.Range("O33").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P33").Value = "early_exit"
    $ws.Range("Q33").Value = 0.13
}
